# "Ajout données sur le toit" - add the roof segment data table to Feuil1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Segment"
$ws.Range("B1").Value = "angle"
$ws.Range("C1").Value = "surface"
$ws.Range("D1").Value = "coordonnées"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = "(x1,y1,z1;…)"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 100

$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# Best-fit the column widths to the new content (Excel auto-sizes columns
# to fit the typed values, so the saved file carries explicit widths for
# columns A:D sized to "Segment"/"angle"/"surface"/"coordonnées" and the
# longest value underneath each). The host stores ColumnWidth + 5/6 (its
# internal glyph-padding constant) as the saved <col width>, so the desired
# on-disk widths are requested net of that fixed offset.
$padding = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 8.85546875 - $padding
$ws.Columns.Item(2).ColumnWidth = 5.85546875 - $padding
$ws.Columns.Item(3).ColumnWidth = 7.42578125 - $padding
$ws.Columns.Item(4).ColumnWidth = 12.5703125 - $padding

# Leave the selection on D3, matching the saved view state.
$ws.Range("D3").Select() | Out-Null
